$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0552
$ws.Range("E2").Value = 0.04395
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 171.78
$ws.Range("L2").Value = 0.3107452966714905
$ws.Range("M2").Value = 52.08
$ws.Range("N2").Value = 0.03054545454545455
$ws.Range("O2").Value = 0.3031784841075795
$ws.Range("P2").Value = 52.08
$ws.Range("Q2").Value = 0.03054545454545455
$ws.Range("R2").Value = 0.3031784841075795
$ws.Range("U2").Value = 1999.7
$ws.Range("V2").Value = 1.172844574780059
$ws.Range("W2").Value = 0.05702970297029702
$ws.Range("X2").Value = 0.06937372615854384
$ws.Range("Y2").Value = -0.01234402318824682
$ws.Range("Z2").Value = 0.1785645067510821
$ws.Range("AB2").Value = 0.04717868226609809
$ws.Range("AC2").Value = -0.04717868226609809
$ws.Range("AD2").Value = 2500
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2500
$ws.Range("AG2").Value = 500.3
$ws.Range("AH2").Value = 0.5945303210463734
$ws.Range("AI2").Value = 0.5129362522825663
$ws.Range("AJ2").Value = 0.2268625583820795
$ws.Range("AK2").Value = 0.174065827012734
$ws.Range("D3").Value = 0.092
$ws.Range("E3").Value = 0.04219999999999999
$ws.Range("K3").Value = 2.88
$ws.Range("L3").Value = 0.2526315789473684
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 52.9
$ws.Range("V3").Value = 1.574404761904762
$ws.Range("W3").Value = 0.05702970297029702
$ws.Range("X3").Value = 0.05194855475759563
$ws.Range("Y3").Value = 0.005081148212701392
$ws.Range("Z3").Value = 0.456
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04716336072554066
$ws.Range("AC3").Value = -0.04716336072554066
$ws.Range("AD3").Value = 16.6
$ws.Range("AF3").Value = 16.6
$ws.Range("AG3").Value = -36.3
$ws.Range("AH3").Value = 0.3306772908366534
$ws.Range("AI3").Value = 0.2597809076682316
$ws.Range("AJ3").Value = 13.44444444444446
$ws.Range("AK3").Value = -3.3
$ws.Range("B4").Value = "MCB Group Limited (MUSE:MCBG.N0000)"
$ws.Range("D4").Value = 0.0552
$ws.Range("E4").Value = 0.0457
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 184.9
$ws.Range("L4").Value = 0.4547466797835711
$ws.Range("M4").Value = 45.6
$ws.Range("N4").Value = 0.03172615320392402
$ws.Range("O4").Value = 0.2466197944835046
$ws.Range("P4").Value = 45.6
$ws.Range("Q4").Value = 0.03172615320392402
$ws.Range("R4").Value = 0.2466197944835046
$ws.Range("U4").Value = 1438.5
$ws.Range("V4").Value = 1.000834898768524
$ws.Range("W4").Value = 0.1187692702980473
$ws.Range("X4").Value = 0.06937372615854384
$ws.Range("Y4").Value = 0.04939554413950344
$ws.Range("Z4").Value = 0.2427607618365276
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04717868226609809
$ws.Range("AC4").Value = -0.04717868226609809
$ws.Range("AD4").Value = 1691.1
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1691.1
$ws.Range("AG4").Value = 252.5999999999999
$ws.Range("AH4").Value = 0.5405638665132336
$ws.Range("AI4").Value = 0.5008292365100989
$ws.Range("AJ4").Value = 0.1494763003728031
$ws.Range("AK4").Value = 0.1303338321036066
$ws.Range("B5").Value = "SBM Holdings Ltd (MUSE:SBMH.N0000)"
$ws.Range("D5").Value = 0.0344
$ws.Range("K5").Value = -16
$ws.Range("L5").Value = -0.1186943620178041
$ws.Range("M5").Value = 6.48
$ws.Range("N5").Value = 0.02768047842802222
$ws.Range("O5").Value = -0.405
$ws.Range("P5").Value = 6.48
$ws.Range("Q5").Value = 0.02768047842802222
$ws.Range("R5").Value = -0.405
$ws.Range("U5").Value = 508.3
$ws.Range("V5").Value = 2.171294318667236
$ws.Range("W5").Value = -0.02224384818573613
$ws.Range("X5").Value = 0.1257409445287726
$ws.Range("Y5").Value = -0.1479847927145087
$ws.Range("Z5").Value = 0.09656852210043701
$ws.Range("AB5").Value = 0.050535527674282
$ws.Range("AC5").Value = -0.050535527674282
$ws.Range("AD5").Value = 792.3
$ws.Range("AF5").Value = 792.3
$ws.Range("AG5").Value = 283.9999999999999
$ws.Range("AH5").Value = 0.7719212782540921
$ws.Range("AI5").Value = 0.5527417329426537
$ws.Range("AJ5").Value = 0.5481567265006755
$ws.Range("AK5").Value = 0.3069938385039455

# Clears (cell removals)
$ws.Range("F2").ClearContents()
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("E5").ClearContents()
